$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.688.59"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.557.38"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.31"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.65"
$ws.Range("E6").Value = "  +7.22%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.04"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +9.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "2.543.28"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.881"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.47"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").Value = "42.790.47"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.47"
$ws.Range("E18").Value = "  +8.07%  "
$ws.Range("D19").Value = "0.0₃0988"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.60"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "256.75"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.95"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "28.14"
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.96"
$ws.Range("E27").Value = "  +8.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.94"
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0804"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.31"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.33"
$ws.Range("E36").Value = "  +14.75%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.70"
$ws.Range("E37").Value = "  +7.45%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.85"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.06"
$ws.Range("E41").Value = "  +30.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.067.74"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.94"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("E47").Value = "  +5.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.14"
$ws.Range("E48").Value = "  +10.93%  "
$ws.Range("D49").Value = "2.803.30"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.62"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  +2.39%  "
